$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add description (D column) to the existing security-group table (rows 10-13)
$ws.Range("D10").Value = "인터넷 - http"
$ws.Range("D11").Value = "인터넷보안 -https"
$ws.Range("D12").Value = "Flask : TCP"
$ws.Range("D13").Value = "Maria / MySQL : TCP "

# New section title
$ws.Range("A15").Value = "보안그룹"

# Re-create the header rows of the table (mirrors rows 8-9)
$ws.Range("A16").Value = "범위"
$ws.Range("B16").Value = "포트"
$ws.Range("D16").Value = "설명"

$ws.Range("B17").Value = "인"
$ws.Range("C17").Value = "아웃"

# Data rows (mirrors rows 10-13, now including the description column)
$ws.Range("A18").Value = "10.0.0.0/22"
$ws.Range("B18").Value = 80
$ws.Range("C18").Value = 80
$ws.Range("D18").Value = "인터넷 - http"

$ws.Range("B19").Value = 443
$ws.Range("C19").Value = 443
$ws.Range("D19").Value = "인터넷보안 -https"

$ws.Range("A20").Value = "10.0.8.0/23"
$ws.Range("B20").Value = 5000
$ws.Range("C20").Value = 5000
$ws.Range("D20").Value = "Flask : TCP"

$ws.Range("A21").Value = "10.0.13.0/24"
$ws.Range("B21").Value = 3306
$ws.Range("C21").Value = 3306
$ws.Range("D21").Value = "Maria / MySQL : TCP "

# Update current selection to match the final workbook state
$ws.Range("H17").Select()
